# Migrate participant storage to Google Sheets
# Append a new participant (IvanZ) as row 3, mirroring the existing
# row-2 record layout: Имя, Username, ChatID, Email, Формат, Оплата, Фидбэк

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry over the blank "Фидбэк" placeholder cell from row 2 so the new
# row keeps the same trailing empty column, then fill in the new data.
$ws.Range("G2").Copy($ws.Range("G3"))

$ws.Range("A3").Value = "IvanZ"
$ws.Range("B3").Value = "@ivan_z89"
$ws.Range("C3").Value = 6172894470
$ws.Range("D3").Value = "qwe@qwe.ty"
$ws.Range("E3").Value = "free"
$ws.Range("F3").Value = "нет"
